$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row: a second set of test credentials (testuser2 / Test@123)
$ws.Range("A3").Value = "testuser2"
$ws.Range("B3").Value = "Test@123"

# Mirror the mailto hyperlink + Hyperlink style that row 2's password cell has
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Test@123")
$ws.Range("B3").Style = $ws.Range("B2").Style

# Keep selection in step with the new last row, like the source workbook
$ws.Range("B4").Select()
